$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "50.962.69"
$ws.Range("E2").Value = "  -1.95%  "
$ws.Range("D3").Value = "2.900.80"
$ws.Range("E3").Value = "  -0.71%  "
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").Value = "346.07"
$ws.Range("E5").Value = "  -1.75%  "
$ws.Range("D6").Value = "106.31"
$ws.Range("E6").Value = "  -5.95%  "
$ws.Range("D7").Value = "0.547"
$ws.Range("E7").Value = "  -2.02%  "
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "0.600"
$ws.Range("E9").Value = "  -2.82%  "
$ws.Range("D10").Value = "37.14"
$ws.Range("E10").Value = "  -5.74%  "
$ws.Range("E11").Value = "  +1.09%  "
$ws.Range("D12").Value = "0.0842"
$ws.Range("E12").Value = "  -3.48%  "
$ws.Range("D13").Value = "18.70"
$ws.Range("E13").Value = "  -6.48%  "
$ws.Range("D14").Value = "3.366.80"
$ws.Range("E14").Value = "  +0.03%  "
$ws.Range("D15").Value = "7.58"
$ws.Range("E15").Value = "  -1.74%  "
$ws.Range("D16").Value = "2.950.32"
$ws.Range("E16").Value = "  +1.63%  "
$ws.Range("D17").Value = "0.950"
$ws.Range("E17").Value = "  -3.35%  "
$ws.Range("D18").Value = "51.029.43"
$ws.Range("E18").Value = "  -1.78%  "
$ws.Range("D19").Value = "3.37"
$ws.Range("E19").Value = "  +2.64%  "
$ws.Range("D20").Value = "7.31"
$ws.Range("E20").Value = "  -3.38%  "
$ws.Range("D21").Value = "13.25"
$ws.Range("E21").Value = "  -5.73%  "
$ws.Range("D22").Value = "0.0₃0951"
$ws.Range("E22").Value = "  -2.50%  "
$ws.Range("D23").Value = "68.34"
$ws.Range("E23").Value = "  -3.69%  "
$ws.Range("D24").Value = "258.86"
$ws.Range("E24").Value = "  -3.50%  "
$ws.Range("D25").Value = "2.66"
$ws.Range("E25").Value = "  -4.19%  "
$ws.Range("D26").Value = "0.169"
$ws.Range("E26").Value = "  -4.88%  "
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("D28").Value = "26.06"
$ws.Range("E28").Value = "  -2.76%  "
$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D29").Value = "7.42"
$ws.Range("E29").Value = "  +7.07%  "
$ws.Range("D30").Value = "0.103"
$ws.Range("E30").Value = "  +0.73%  "
$ws.Range("D31").Value = "10.08"
$ws.Range("E31").Value = "  -5.05%  "
$ws.Range("D32").Value = "5.98"
$ws.Range("E32").Value = "  -0.22%  "
$ws.Range("D33").Value = "34.92"
$ws.Range("E33").Value = "  -5.44%  "
$ws.Range("D34").Value = "2.12"
$ws.Range("E34").Value = "  +1.26%  "
$ws.Range("D35").Value = "50.42"
$ws.Range("E35").Value = "  -4.92%  "
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.40%  "
$ws.Range("D37").Value = "0.0418"
$ws.Range("E37").Value = "  -7.61%  "
$ws.Range("D38").Value = "3.07"
$ws.Range("E38").Value = "  -8.53%  "
$ws.Range("D39").Value = "17.36"
$ws.Range("E39").Value = "  -6.60%  "
$ws.Range("D40").Value = "1.91"
$ws.Range("E40").Value = "  -6.16%  "
$ws.Range("D41").Value = "2.59"
$ws.Range("E41").Value = "  -3.09%  "
$ws.Range("D42").Value = "0.114"
$ws.Range("E42").Value = "  -1.82%  "
$ws.Range("D43").Value = "21.98"
$ws.Range("E43").Value = "  -3.51%  "
$ws.Range("E44").Value = "  +6.92%  "
$ws.Range("E45").Value = "  -2.49%  "
$ws.Range("D46").Value = "2.075.08"
$ws.Range("E46").Value = "  -5.08%  "
$ws.Range("D47").Value = "3.26"
$ws.Range("E47").Value = "  -6.55%  "
$ws.Range("D48").Value = "2.24"
$ws.Range("E48").Value = "  -10.80%  "
$ws.Range("D49").Value = "0.236"
$ws.Range("E49").Value = "  -4.99%  "
$ws.Range("D50").Value = "0.0330"
$ws.Range("E50").Value = "  -3.00%  "
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").Value = "59.66"
$ws.Range("E51").Value = "  -4.69%  "
